$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add 13 new data rows (417-429) below the existing data block, which
# currently ends at row 416. Each new row records the same production
# line ("ASSY COVER-OVER HEAD CONSOLE") for date serial 46066
# (2026-02-13), mirroring the layout/format of the preceding rows.
# ------------------------------------------------------------------

$newRows = 417..429
$sourceRow = 416

foreach ($r in $newRows) {
    # Copy the cell formatting (number format / style) from the last
    # existing row (416, columns B:F) down onto this new row so the new
    # cells pick up the same look (date style, text style, number style).
    $ws.Range("B$sourceRow").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)

    $ws.Range("C$sourceRow").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)

    $ws.Range("D$sourceRow").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)

    $ws.Range("E$sourceRow").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)

    $ws.Range("F$sourceRow").Copy()
    $ws.Range("F$r").PasteSpecial(-4122)

    # Column G on these new rows carries the same (non-standard) number
    # style as D:F rather than the usual G-column style, so pull the
    # format from F instead of the existing G column.
    $ws.Range("F$sourceRow").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)

    # Now write the actual values for this row.
    $ws.Range("B$r").Value = 46066
    $ws.Range("C$r").Value = "ASSY COVER-OVER HEAD CONSOLE"
    $ws.Range("D$r").Value = 1112
    $ws.Range("E$r").Value = 1112
    $ws.Range("F$r").Value = 112
    $ws.Range("G$r").Value = 112
}

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Update the sheet's current selection to match where the user ended up
# after entering the new rows.
# ------------------------------------------------------------------
[void]$ws.Range("G417:G429").Select()
